$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# NOTE: Shape.Left/.Top/.Width/.Height (and AddTextbox args) are in POINTS,
# like real PowerPoint COM automation (1 pt = 12700 EMU).

# --- Move the Subtitle placeholder up ---
$subtitle = $s.Shapes.Item("Subtitle 2")
$subtitle.Left = 60
$subtitle.Top = 348

# --- Move the connector line up to match ---
$connector = $s.Shapes.Item("Straight Connector 7")
$connector.Left = 672
$connector.Top = 360

# --- Add a textbox with the event name/date at the bottom of the slide ---
$textbox = $s.Shapes.AddTextbox(1, 6, 504.91874, 702, 29.08126)
$textbox.Name = "TextBox 3"

$tf = $textbox.TextFrame
$tf.WordWrap = -1
$tf.AutoSize = 1

$tr = $tf.TextRange
$tr.Text = "ML Workshop (Copenhagen, 13 September, 2012)"
$tr.Font.Bold = -1
$tr.Font.Name = "+mj-lt"

$run1 = $tr.Characters(1, 12)
$run1.Font.Bold = -1

$run2 = $tr.Characters(13, 33)
$run2.Font.Bold = 0
